$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the summary header values (VALOR MORA total, worker/period counts)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 12418
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2

# ---------------------------------------------------------------------------
# 2. Remove the old detail rows that belonged to the workers/periods that are
#    no longer part of this statement (the block used to run from row 19
#    through row 29 - this also slides the old "last" detail row (30) plus
#    the signature rows up into their new positions).
# ---------------------------------------------------------------------------
$ws.Rows("19:29").Delete()

# ---------------------------------------------------------------------------
# 3. Rewrite the remaining detail rows (16-19) with the new account data.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "73581195"
$ws.Range("D16").Value = "EDER LUIS MENDEZ MARTINEZ"
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 4417
$ws.Range("G16").Value = 908526

$ws.Range("C17").Value = "1002322475"
$ws.Range("D17").Value = "JHON JAIRO PALOMINO TERAN"
$ws.Range("E17").Value = "2204"
$ws.Range("F17").Value = 2667
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "1007188262"
$ws.Range("D18").Value = "DEIBER AVILA TERAN"
$ws.Range("E18").Value = "2204"
$ws.Range("F18").Value = 2667
$ws.Range("G18").Value = 1000000

$ws.Range("C19").Value = "1007315342"
$ws.Range("D19").Value = "JOYNER LUIS ROMERO RODRIGUEZ"
$ws.Range("E19").Value = "2204"
$ws.Range("F19").Value = 2667
$ws.Range("G19").Value = 1000000
